$wb = $excel.ActiveWorkbook

# --- ALERTS ---
$ws_ALERTS = $wb.Worksheets.Item("ALERTS")
$ws_ALERTSData = @(
    @(14, '2026-01-30', '14:36:18', '14:00', 'Living Room', 'CRITICAL EMERGENCY', 'FALL_DETECTED'),
    @(15, '2026-01-30', '14:36:21', '14:00', 'Living Room', 'CRITICAL EMERGENCY', 'FALL_DETECTED'),
    @(16, '2026-01-30', '14:36:59', '14:00', 'Living Room', 'CRITICAL EMERGENCY', 'FALL_DETECTED'),
    @(17, '2026-01-30', '14:40:20', '14:00', 'Living Room', 'CRITICAL EMERGENCY', 'FALL_DETECTED'),
    @(18, '2026-01-30', '14:40:51', '14:00', 'Living Room', 'CRITICAL EMERGENCY', 'FALL_DETECTED'),
    @(19, '2026-01-30', '14:41:54', '14:00', 'Living Room', 'CRITICAL EMERGENCY', 'FALL_DETECTED'),
    @(20, '2026-01-30', '14:41:58', '14:00', 'Living Room', 'CRITICAL EMERGENCY', 'FALL_DETECTED'),
)
foreach ($row in $ws_ALERTSData) {
    $r = $row[0]
    $ws_ALERTS.Cells.Item($r, 1).NumberFormat = "@"
    $ws_ALERTS.Cells.Item($r, 1).Value = $row[1]
    $ws_ALERTS.Cells.Item($r, 2).Value = $row[2]
    $ws_ALERTS.Cells.Item($r, 3).Value = $row[3]
    $ws_ALERTS.Cells.Item($r, 4).Value = $row[4]
    $ws_ALERTS.Cells.Item($r, 5).Value = $row[5]
    $ws_ALERTS.Cells.Item($r, 6).Value = $row[6]
}

# --- PIR ---
$ws_PIR = $wb.Worksheets.Item("PIR")
$ws_PIRData = @(
    @(128, '2026-01-30', '14:39:04', '14:00', 'Living Room', 'RECOVERY_DETECTION', 'Inactive'),
    @(129, '2026-01-30', '14:40:23', '14:00', 'Living Room', 'RECOVERY_DETECTION', 'Inactive'),
)
foreach ($row in $ws_PIRData) {
    $r = $row[0]
    $ws_PIR.Cells.Item($r, 1).NumberFormat = "@"
    $ws_PIR.Cells.Item($r, 1).Value = $row[1]
    $ws_PIR.Cells.Item($r, 2).Value = $row[2]
    $ws_PIR.Cells.Item($r, 3).Value = $row[3]
    $ws_PIR.Cells.Item($r, 4).Value = $row[4]
    $ws_PIR.Cells.Item($r, 5).Value = $row[5]
    $ws_PIR.Cells.Item($r, 6).Value = $row[6]
}

# --- Proximity ---
$ws_Proximity = $wb.Worksheets.Item("Proximity")
$ws_ProximityData = @(
    @(64, '2026-01-30', '14:32:51', '14:00', 'Living Room Main Door', 'ENTER', 'User ENTERED Living Room Main Door'),
    @(65, '2026-01-30', '14:32:55', '14:00', 'Living Room Main Door', 'EXIT', 'User EXITED Living Room Main Door'),
    @(66, '2026-01-30', '14:33:02', '14:00', 'Living Room Main Door', 'ENTER', 'User ENTERED Living Room Main Door'),
    @(67, '2026-01-30', '14:33:05', '14:00', 'Living Room Main Door', 'EXIT', 'User EXITED Living Room Main Door'),
)
foreach ($row in $ws_ProximityData) {
    $r = $row[0]
    $ws_Proximity.Cells.Item($r, 1).NumberFormat = "@"
    $ws_Proximity.Cells.Item($r, 1).Value = $row[1]
    $ws_Proximity.Cells.Item($r, 2).Value = $row[2]
    $ws_Proximity.Cells.Item($r, 3).Value = $row[3]
    $ws_Proximity.Cells.Item($r, 4).Value = $row[4]
    $ws_Proximity.Cells.Item($r, 5).Value = $row[5]
    $ws_Proximity.Cells.Item($r, 6).Value = $row[6]
}

# --- mmWave ---
$ws_mmWave = $wb.Worksheets.Item("mmWave")
$ws_mmWaveData = @(
    @(82, '2026-01-30', '14:32:48', '14:00', 'Living Room', 'FALL_DETECTED', 'EMERGENCY'),
    @(83, '2026-01-30', '14:32:49', '14:00', 'Living Room', 'FALL_DETECTED', 'EMERGENCY'),
    @(84, '2026-01-30', '14:36:18', '14:00', 'Living Room', 'CRITICAL EMERGENCY', 'FALL_DETECTED'),
    @(85, '2026-01-30', '14:36:21', '14:00', 'Living Room', 'CRITICAL EMERGENCY', 'FALL_DETECTED'),
    @(86, '2026-01-30', '14:36:26', '14:00', 'Living Room', 'PRESENCE_DETECTED', 'Active'),
    @(87, '2026-01-30', '14:36:36', '14:00', 'Living Room', 'PRESENCE_DETECTED', 'Active'),
    @(88, '2026-01-30', '14:36:59', '14:00', 'Living Room', 'CRITICAL EMERGENCY', 'FALL_DETECTED'),
    @(89, '2026-01-30', '14:39:04', '14:00', 'Living Room', 'PRESENCE_DETECTED', 'Active'),
    @(90, '2026-01-30', '14:39:14', '14:00', 'Living Room', 'PRESENCE_DETECTED', 'Active'),
    @(91, '2026-01-30', '14:39:25', '14:00', 'Living Room', 'PRESENCE_DETECTED', 'Active'),
    @(92, '2026-01-30', '14:39:35', '14:00', 'Living Room', 'PRESENCE_DETECTED', 'Active'),
    @(93, '2026-01-30', '14:39:46', '14:00', 'Living Room', 'PRESENCE_DETECTED', 'Active'),
    @(94, '2026-01-30', '14:39:56', '14:00', 'Living Room', 'PRESENCE_DETECTED', 'Active'),
    @(95, '2026-01-30', '14:40:20', '14:00', 'Living Room', 'CRITICAL EMERGENCY', 'FALL_DETECTED'),
    @(96, '2026-01-30', '14:40:23', '14:00', 'Living Room', 'PRESENCE_DETECTED', 'Active'),
    @(97, '2026-01-30', '14:40:51', '14:00', 'Living Room', 'CRITICAL EMERGENCY', 'FALL_DETECTED'),
    @(98, '2026-01-30', '14:41:54', '14:00', 'Living Room', 'CRITICAL EMERGENCY', 'FALL_DETECTED'),
    @(99, '2026-01-30', '14:41:58', '14:00', 'Living Room', 'CRITICAL EMERGENCY', 'FALL_DETECTED'),
    @(100, '2026-01-30', '14:42:10', '14:00', 'Living Room', 'NO_MOTION_DETECTED', 'Inactive'),
    @(101, '2026-01-30', '14:42:21', '14:00', 'Living Room', 'PRESENCE_DETECTED', 'Active'),
    @(102, '2026-01-30', '14:42:31', '14:00', 'Living Room', 'PRESENCE_DETECTED', 'Active'),
    @(103, '2026-01-30', '14:42:50', '14:00', 'Living Room', 'CRITICAL EMERGENCY', 'FALL_DETECTED'),
    @(104, '2026-01-30', '14:42:54', '14:00', 'Living Room', 'CRITICAL EMERGENCY', 'FALL_DETECTED'),
    @(105, '2026-01-30', '14:42:54', '14:00', 'Living Room', 'PRESENCE_DETECTED', 'Active'),
    @(106, '2026-01-30', '14:43:03', '14:00', 'Living Room', 'PRESENCE_DETECTED', 'Active'),
    @(107, '2026-01-30', '14:43:13', '14:00', 'Living Room', 'PRESENCE_DETECTED', 'Active'),
    @(108, '2026-01-30', '14:43:24', '14:00', 'Living Room', 'PRESENCE_DETECTED', 'Active'),
    @(109, '2026-01-30', '14:43:34', '14:00', 'Living Room', 'PRESENCE_DETECTED', 'Active'),
)
foreach ($row in $ws_mmWaveData) {
    $r = $row[0]
    $ws_mmWave.Cells.Item($r, 1).NumberFormat = "@"
    $ws_mmWave.Cells.Item($r, 1).Value = $row[1]
    $ws_mmWave.Cells.Item($r, 2).Value = $row[2]
    $ws_mmWave.Cells.Item($r, 3).Value = $row[3]
    $ws_mmWave.Cells.Item($r, 4).Value = $row[4]
    $ws_mmWave.Cells.Item($r, 5).Value = $row[5]
    $ws_mmWave.Cells.Item($r, 6).Value = $row[6]
}

# --- Camera ---
$ws_Camera = $wb.Worksheets.Item("Camera")
$ws_CameraData = @(
    @(9, '2026-01-30', '14:32:51', '14:00', 'Living Room Main Door', 'Image Captured', 'Active'),
    @(10, '2026-01-30', '14:33:03', '14:00', 'Living Room Main Door', 'Image Captured', 'Active'),
)
foreach ($row in $ws_CameraData) {
    $r = $row[0]
    $ws_Camera.Cells.Item($r, 1).NumberFormat = "@"
    $ws_Camera.Cells.Item($r, 1).Value = $row[1]
    $ws_Camera.Cells.Item($r, 2).Value = $row[2]
    $ws_Camera.Cells.Item($r, 3).Value = $row[3]
    $ws_Camera.Cells.Item($r, 4).Value = $row[4]
    $ws_Camera.Cells.Item($r, 5).Value = $row[5]
    $ws_Camera.Cells.Item($r, 6).Value = $row[6]
}
